# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets contain identical data, and the same set of rows changed.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3033
    6  = 18
    7  = 1656
    9  = 83
    11 = 1361
    12 = 11
    13 = 497
    14 = 346
    15 = 24
    16 = 72
    18 = 122
    20 = 105
    21 = 3157
    22 = 387
    23 = 120
    24 = 201
    25 = 6
    27 = 91
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
